# "feat: change framework route manger"
#
# The enum sheet (全局枚举表) gets a brand-new, empty column inserted right
# before column C — every existing column from C onward (old C..J) shifts
# one slot to the right (new D..K). The inserted column inherits the
# width of column B (18.5 "chars"), and the view/selection state left
# behind by the editing session is updated on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 全局枚举表
$ws2 = $wb.Worksheets.Item(2)   # 生成表

# Insert a new blank column at C; B..old-C.. shift right by one.
$ws1.Columns.Item(3).Insert()

# New column C should look like column B (same width).
$ws1.Columns.Item(3).ColumnWidth = 18.5 - (5/7)

# Restore the view/selection state recorded at save time.
# (select on sheet2 first so sheet1 ends up as the active/selected tab)
$ws2.Range("A1").Select()
$ws1.Range("B18").Select()
